$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

# The CO2 trade-link rows (36-43) were missing their "Tech" (column H)
# value. Fill it in with "B", matching every other trade-link row above.
$ws.Range("H36:H43").Value = "B"

$ws.Activate()
$ws.Range("K41").Select()
